# Apply "average with safety stocks" edit.
# 1) Productdata sheet: rescale SetupCosts (D), BackorderCosts (F) and
#    LostSale (I) columns for rows 2-11 by a factor of 0.0004.
# 2) ForcastedStandardDeviation sheet: zero out the standard-deviation
#    values for buckets 7-9 (rows 9-11), columns B-E.

$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("Productdata")

$productUpdates = @{
    2  = @{ D = 0.0016; F = 0.016; I = 0.16 }
    3  = @{ D = 0.0028; F = 0.028; I = 0.28 }
    4  = @{ D = 0.0024; F = 0.024; I = 0.24 }
    5  = @{ D = 0.0012; F = 0.012; I = 0.12 }
    6  = @{ D = 0.0012; F = 0.012; I = 0.12 }
    7  = @{ D = 0.0012; F = 0.012; I = 0.12 }
    8  = @{ D = 0.0008; F = 0.008; I = 0.08 }
    9  = @{ D = 0.0004; F = 0.004; I = 0.04 }
    10 = @{ D = 0.0004; F = 0.004; I = 0.04 }
    11 = @{ D = 0.0004; F = 0.004; I = 0.04 }
}

foreach ($row in $productUpdates.Keys) {
    $vals = $productUpdates[$row]
    $wsProduct.Range("D$row").Value = $vals.D
    $wsProduct.Range("F$row").Value = $vals.F
    $wsProduct.Range("I$row").Value = $vals.I
}

$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")

foreach ($row in 9..11) {
    $wsStdDev.Range("B$row`:E$row").Value = 0
}
